$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.Value = "'" + $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "293.70"
Set-TextValue "E2" "1.35%"
Set-TextValue "D3" "31.06"
Set-TextValue "E3" "0.59%"
Set-TextValue "D4" "4.935"
Set-TextValue "E4" "1.06%"
Set-TextValue "D5" "0.07350"
Set-TextValue "E5" "2.76%"
Set-TextValue "D6" "2.274"
Set-TextValue "E6" "23.55%"
Set-TextValue "D7" "7.697"
Set-TextValue "E7" "0.78%"
Set-TextValue "D8" "3.786"
Set-TextValue "E8" "0.49%"
Set-TextValue "D9" "0.9121"
Set-TextValue "E9" "1.97%"
Set-TextValue "D10" "0.1685"
Set-TextValue "E10" "2.75%"
Set-TextValue "D11" "0.08166"
Set-TextValue "E11" "8.40%"
Set-TextValue "D12" "0.08268"
Set-TextValue "E12" "2.11%"
Set-TextValue "D13" "0.03104"
Set-TextValue "E13" "3.78%"
Set-TextValue "E14" "0.70%"
Set-TextValue "D15" "0.001527"
Set-TextValue "E15" "2.11%"
Set-TextValue "D16" "0.005706"
Set-TextValue "E16" "-1.10%"
Set-TextValue "E17" "0.68%"
Set-TextValue "D18" "2.081"
Set-TextValue "E18" "-1.49%"
Set-TextValue "D19" "0.3329"
Set-TextValue "E19" "1.56%"
Set-TextValue "E20" "0.84%"
Set-TextValue "D21" "3.979"
Set-TextValue "E21" "-6.69%"
Set-TextValue "D22" "0.2101"
Set-TextValue "E22" "5.00%"
Set-TextValue "D23" "0.04544"
Set-TextValue "E23" "1.62%"
Set-TextValue "E24" "-0.06%"
Set-TextValue "D25" "0.004341"
Set-TextValue "E25" "-6.81%"
Set-TextValue "E26" "4.05%"
Set-TextValue "D27" "0.0003394"
Set-TextValue "D39" "0.01604"
Set-TextValue "E39" "-2.00%"
Set-TextValue "D40" "0.04433"
Set-TextValue "E40" "2.05%"
Set-TextValue "D41" "0.007354"
Set-TextValue "E41" "0.01%"
Set-TextValue "D42" "0.008747"
Set-TextValue "E43" "1.60%"
Set-TextValue "D44" "0.002111"
Set-TextValue "E44" "5.37%"
Set-TextValue "D45" "0.009202"
Set-TextValue "E45" "-10.07%"
Set-TextValue "D46" "0.00005948"
Set-TextValue "E46" "1.48%"
Set-TextValue "E47" "-0.12%"
Set-TextValue "E48" "1.51%"
Set-TextValue "D50" "0.00002100"
Set-TextValue "E50" "-0.12%"
Set-TextValue "D51" "0.0002000"
Set-TextValue "E51" "-0.12%"
